# Logged Week 16 and performed season sim from Week 17
# Update the "Road" (R) row totals on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 3 ("R") ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 431
$wsOff.Range("C3").Value = 302
$wsOff.Range("D3").Value = 115
$wsOff.Range("E3").Value = 57
$wsOff.Range("F3").Value = 10
$wsOff.Range("G3").Value = 6

# --- DEF sheet: row 3 ("R") ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 477
$wsDef.Range("C3").Value = 355
$wsDef.Range("D3").Value = 108
$wsDef.Range("E3").Value = 54
$wsDef.Range("F3").Value = 9
